$d = $word.ActiveDocument

function Get-ParaRange($idx) {
    $pp = $d.Paragraphs($idx)
    return $d.Range($pp.Range.Start, $pp.Range.End - 1)
}

function Get-ParaEnd($idx) {
    return (Get-ParaRange $idx).End
}

# Appends $text (as its own new run, formatted italic to match the
# surrounding "ListParagraph" runs) to the end of paragraph $paraIdx,
# without letting it merge into the immediately preceding run.
function Append-Segment($paraIdx, $text) {
    $endPos = Get-ParaEnd $paraIdx
    $markerRng = $d.Range($endPos, $endPos)
    $markerRng.InsertAfter("X")
    $endPos2 = Get-ParaEnd $paraIdx
    $bmRng = $d.Range($endPos2 - 1, $endPos2)
    $d.Bookmarks.Add("TMPSEG", $bmRng)
    $bm = $d.Bookmarks("TMPSEG")
    $bmRange = $bm.Range
    $bmRange.Text = $text
    $bmRange.Italic = 1
    $bmRange.ItalicBi = 1
    $d.Bookmarks("TMPSEG").Delete()
}

# Locate the "Số repository có: 14" paragraph (the filled-in first
# collaborator block) so we can find the (currently empty) second
# collaborator block that immediately follows it.
$srcIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt.Contains("repository") -and $ptxt.Contains("14")) {
        $srcIdx = $i
        break
    }
}

$nameIdx  = $srcIdx + 1   # "Tên đầy đủ (người thứ nhất): "
$linkIdx  = $srcIdx + 2   # "Link tài khoản Github : "
$repoIdx  = $srcIdx + 3   # "Số repository có: "

# 1) Move the "_GoBack" bookmark off the "...: 14" paragraph - it will be
#    re-created at the end of the newly filled-in repo-count paragraph.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2) "Tên đầy đủ (người thứ nhất): " -> split into the six runs the
#    author's edit produced: "...thứ" / " hai" / "):" / " " / " " / "Isaac Flores"
$nameRng = Get-ParaRange $nameIdx
$nameRng.Text = "Tên đầy đủ (người thứ"
Append-Segment $nameIdx " hai"
Append-Segment $nameIdx "):"
Append-Segment $nameIdx " "
Append-Segment $nameIdx " "
Append-Segment $nameIdx "Isaac Flores"

# 3) "Link tài khoản Github : " -> append the Github URL as a new run.
Append-Segment $linkIdx "https://github.com/isaacaflores2"

# 4) "Số repository có: " -> append "13" as a new run, then re-create the
#    "_GoBack" bookmark immediately around that new run.
Append-Segment $repoIdx "13"

$repoEnd = Get-ParaEnd $repoIdx
$bmFinalRng = $d.Range($repoEnd - 2, $repoEnd)
$d.Bookmarks.Add("_GoBack", $bmFinalRng)
